# Auto-generated Excel COM-interop script to apply the numeric updates
# described in the commit diff across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 202.4
$ws.Range("I28").Value = 200.71428
$ws.Range("J28").Value = 206.33333
$ws.Range("K28").Value = 200.71428
$ws.Range("L28").Value = 206.33333
$ws.Range("M28").Value = 284.28572
$ws.Range("N28").Value = -1176.33333

$ws.Range("H40").Value = 2320.9285
$ws.Range("J40").Value = 1653.8182
$ws.Range("L40").Value = 1653.8182
$ws.Range("N40").Value = -2003.8182

$ws.Range("H113").Value = 2818.077
$ws.Range("J113").Value = 3241.6667
$ws.Range("L113").Value = 3241.6667
$ws.Range("N113").Value = -9749.6667

$ws.Range("H137").Value = 1219.0333
$ws.Range("I137").Value = 1023.8333
$ws.Range("K137").Value = 3071.4999
$ws.Range("M137").Value = -521.4998999999998

$ws.Range("H138").Value = 631832.4399999999
$ws.Range("J138").Value = 871315.2
$ws.Range("L138").Value = 2613945.6
$ws.Range("N138").Value = -2624225.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2936.4329
$ws.Range("I32").Value = 2933.6724
$ws.Range("J32").Value = 2954.2222
$ws.Range("K32").Value = 2933.6724
$ws.Range("L32").Value = 2954.2222
$ws.Range("M32").Value = -2646.6724
$ws.Range("N32").Value = -3528.2222

$ws.Range("H74").Value = 1529.52
$ws.Range("I74").Value = 649.9286
$ws.Range("K74").Value = 649.9286
$ws.Range("M74").Value = 224.0714

$ws.Range("H77").Value = 1529.52
$ws.Range("I77").Value = 649.9286
$ws.Range("K77").Value = 3249.643
$ws.Range("M77").Value = 1118.357

$ws.Range("H102").Value = 13891892
$ws.Range("I102").Value = 15154427
$ws.Range("K102").Value = 15154427
$ws.Range("M102").Value = -15152805

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H87").Value = 75000
$ws.Range("J87").Value = 75000
$ws.Range("L87").Value = 75000
$ws.Range("N87").Value = -77496

$ws.Range("H88").Value = 34000
$ws.Range("J88").Value = 34000
$ws.Range("L88").Value = 34000
$ws.Range("N88").Value = -34812

$ws.Range("H90").Value = 75000
$ws.Range("J90").Value = 75000
$ws.Range("L90").Value = 225000
$ws.Range("N90").Value = -237480

$ws.Range("H91").Value = 34000
$ws.Range("J91").Value = 34000
$ws.Range("L91").Value = 34000
$ws.Range("N91").Value = -36808

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1947.6316
$ws.Range("I31").Value = 1653.3334
$ws.Range("K31").Value = 1653.3334
$ws.Range("M31").Value = -1358.3334

$ws.Range("H34").Value = 1947.6316
$ws.Range("I34").Value = 1653.3334
$ws.Range("K34").Value = 1653.3334
$ws.Range("M34").Value = -1451.3334

$ws.Range("H58").Value = 1496.0667
$ws.Range("J58").Value = 1999
$ws.Range("L58").Value = 1999
$ws.Range("N58").Value = -2405

$ws.Range("H107").Value = 710.1667
$ws.Range("I107").Value = 503.66666
$ws.Range("K107").Value = 503.66666
$ws.Range("M107").Value = 1416.33334

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H125").Value = 17983.334
$ws.Range("J125").Value = 17983.334
$ws.Range("L125").Value = 17983.334
$ws.Range("N125").Value = -22903.334

$ws.Range("H136").Value = 1496.0667
$ws.Range("J136").Value = 1999
$ws.Range("L136").Value = 5997
$ws.Range("N136").Value = -11097

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H96").Value = 8200
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 8200
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 24600
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -28718

$ws.Range("H122").Value = 889.5454999999999
$ws.Range("I122").Value = 600
$ws.Range("J122").Value = 1237
$ws.Range("K122").Value = 5400
$ws.Range("L122").Value = 11133
$ws.Range("M122").Value = -2950
$ws.Range("N122").Value = -16033

$ws.Range("H123").Value = 2525.2
$ws.Range("I123").Value = 1030
$ws.Range("J123").Value = 2899
$ws.Range("K123").Value = 3090
$ws.Range("L123").Value = 8697
$ws.Range("M123").Value = -640
$ws.Range("N123").Value = -13597

$ws.Range("H133").Value = 4155.5
$ws.Range("I133").Value = 2138
$ws.Range("J133").Value = 4705.727
$ws.Range("K133").Value = 6414
$ws.Range("L133").Value = 14117.181
$ws.Range("M133").Value = -1354
$ws.Range("N133").Value = -24237.181

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

$ws.Range("H122").Value = 1406.8276
$ws.Range("I122").Value = 1546.8096
$ws.Range("J122").Value = 1039.375
$ws.Range("K122").Value = 4640.4288
$ws.Range("L122").Value = 3118.125
$ws.Range("M122").Value = -2190.4288
$ws.Range("N122").Value = -8018.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 948
$ws.Range("I22").Value = 516
$ws.Range("J22").Value = 1434
$ws.Range("K22").Value = 516
$ws.Range("L22").Value = 1434
$ws.Range("M22").Value = -221
$ws.Range("N22").Value = -2024

$ws.Range("H27").Value = 948
$ws.Range("I27").Value = 516
$ws.Range("J27").Value = 1434
$ws.Range("K27").Value = 516
$ws.Range("L27").Value = 1434
$ws.Range("M27").Value = -409
$ws.Range("N27").Value = -1648

$ws.Range("H92").Value = 15000
$ws.Range("J92").Value = 15000
$ws.Range("L92").Value = 15000
$ws.Range("N92").Value = -19992

$ws.Range("H110").Value = 30322
$ws.Range("I110").Value = 30000
$ws.Range("J110").Value = 30644
$ws.Range("K110").Value = 30000
$ws.Range("L110").Value = 30644
$ws.Range("M110").Value = -25910
$ws.Range("N110").Value = -38824

$ws.Range("H132").Value = 19153.088
$ws.Range("I132").Value = 1204.5405
$ws.Range("J132").Value = 52357.9
$ws.Range("K132").Value = 3613.6215
$ws.Range("L132").Value = 157073.7
$ws.Range("M132").Value = -1083.6215
$ws.Range("N132").Value = -162133.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 462.2353
$ws.Range("I107").Value = 497.2
$ws.Range("J107").Value = 412.2857
$ws.Range("K107").Value = 1491.6
$ws.Range("L107").Value = 1236.8571
$ws.Range("M107").Value = 428.4000000000001
$ws.Range("N107").Value = -5076.8571

$ws.Range("H119").Value = 19973.5
$ws.Range("J119").Value = 19973.5
$ws.Range("L119").Value = 19973.5
$ws.Range("N119").Value = -29649.5

$ws.Range("H122").Value = 21668788
$ws.Range("J122").Value = 1600
$ws.Range("L122").Value = 4800
$ws.Range("N122").Value = -9700
